$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $rng = $d.Content
    $found = $rng.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        Write-Output "NOT FOUND: $old"
    }
}

# --- Title ---
Replace-Text "Cultural Evolution: The Tapestry of Human Societies" "The Vital Role of History in a Rapidly Changing World"

# --- Author name ---
Replace-Text "Isabella Perez" "Benjamin Walker"

# --- Email address (local-part, then domain; trailing ".edu" run removed) ---
Replace-Text "isabella" "walkerbenjamin24@gmail"
Replace-Text "perez@academicscope" "com"
$rngEdu = $d.Content
$foundEdu = $rngEdu.Find.Execute(".edu", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($foundEdu) {
    $rngEdu.Delete()
}

# --- Body paragraph 1 ---
Replace-Text "The tapestry of human societies is a vibrant and intricate creation, a kaleidoscope of cultures that weaves together beliefs, traditions, and practices" "History is the mosaic of human experiences, an intricate tapestry woven from countless individual journeys"

Replace-Text " Each thread within this tapestry holds a unique story, reflecting the collective aspirations and shared experiences of a people's journey through time" " It paints a vivid portrait of the past, offering invaluable lessons for the present and future"

Replace-Text " Cultural evolution is the dynamic and ever-changing process that shapes these threads, driven by forces of migration, adaptation, and innovation" " In a world marked by unprecedented transformation and uncertainties, understanding history becomes paramount. It enables us to navigate the complex challenges of today and shape a path toward a more promising tomorrow"

Replace-Text "At its core, cultural evolution is an engine of adaptation" "History affords us a kaleidoscope of perspectives, unveiling the intricate interplay of human motivations, choices, and consequences"

Replace-Text " As societies encounter new challenges in their environments, they adapt their beliefs, practices, and technologies" " By retracing the footsteps of those who came before us, we gain insights into our own struggles and aspirations"

Replace-Text " Culture is a living thing, constantly evolving to better suit a society's needs and circumstances" " History provides context, enabling us to grasp the roots of contemporary issues and comprehend the forces that have shaped our present landscape"

Replace-Text " As environmental conditions change - be it climate, disease, or the arrival of new peoples - cultures shift and reshape themselves accordingly, ensuring the continuity of the social fabric" " This understanding fosters empathy, cultivates responsible citizenship, and equips us to address the multifaceted challenges of an ever-changing world"

Replace-Text "Cultural evolution also arises from human interaction" "Moreover, history is a catalyst for critical thinking and innovation"

Replace-Text " When societies encounter one another, whether throughMao Yi , war, or migration, cultural exchange and assimilation become powerful catalysts of change" " It challenges us to examine assumptions, question conventional wisdom, and seek creative solutions to pressing problems"

Replace-Text " New ideas, technologies, and customs flow between cultures, enriching and altering them" " By studying the successes and failures of past societies, we glean valuable lessons that can inform policy decisions, technological advancements, and societal progress"

Replace-Text " The encounters between different ways of life challenge old assumptions and beliefs, leading to the emergence of new perspectives and traditions" " In this way, history serves as a vital compass, orienting us amid the crosscurrents of change and guiding us toward a brighter future"

# --- Summary paragraph ---
Replace-Text "Cultural evolution is an intricate tapestry, woven by the threads of history, adaptation, and interaction" "In a rapidly transforming world, history stands as an indispensable resource for understanding our past, present, and future"

Replace-Text " It is a dynamic process that shapes the beliefs, traditions, and practices of human societies" " It reveals the tapestry of human experiences, offering valuable lessons for navigating the challenges of today"

Replace-Text " Driven by forces such as migration, environmental change, and cultural exchange, cultural evolution helps societies adapt, innovate, and forge new paths in their journeys through time" " History cultivates empathy, fosters responsible citizenship, and ignites critical thinking"

Replace-Text " It is through this dynamic process that the rich tapestry of human culture continues to evolve and thrive, making it an integral and ever-changing part of the human experience" " By examining the triumphs and trials of those who came before us, we gain insights into our own struggles and aspirations, empowering us to shape a future that resonates with lessons learned from the past"

# --- Append a new empty paragraph at the end of the document body ---
$endRange = $d.Content
$endRange.Collapse(0)
$endRange.InsertParagraphAfter()
